$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.730049729347229
$ws.Range("B1").Value = 1.927272200584412
$ws.Range("C1").Value = 2.575516700744629
$ws.Range("D1").Value = 0.9227412939071655
$ws.Range("E1").Value = 1.041294932365417
